# ============================================================================
# Adiciona 2 novas abas: "Ata de Reuniões" (antes de "PROCESSO PADRÃO") e
# "Controle do Projeto" (depois de "PROCESSO PADRÃO").
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Reorganiza as abas
# ---------------------------------------------------------------------------
$procSheet = $wb.Worksheets.Item("PROCESSO PADRÃO")
$ata = $wb.Worksheets.Add($procSheet)
$ata.Name = "Ata de Reuniões"

$procSheet2 = $wb.Worksheets.Item("PROCESSO PADRÃO")
$controle = $wb.Worksheets.Add($null, $procSheet2)
$controle.Name = "Controle do Projeto"

# ---------------------------------------------------------------------------
# 2) Aba "Ata de Reuniões"
# ---------------------------------------------------------------------------
$ws = $ata

$headers = @("Sl#","Data da Reunião","Participanates","Agenda Geral/Topico","Duração`n(hrs)","Action Items","Resp","Data Prevista","Status","Comentários/Notas")
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $headers[$c-1]
}
$headerRange = $ws.Range("A1:J1")
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 9
$headerRange.Interior.Color = 10079487
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$ws.Range("A1:J1").Borders.Item(8).Weight = -4138
$ws.Range("A1:J1").Borders.Item(9).Weight = -4138
$ws.Range("A1:A1").Borders.Item(7).Weight = -4138
$ws.Range("J1:J1").Borders.Item(10).Weight = -4138
$ws.Rows.Item(1).RowHeight = 24

$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "27/10/2016"
$ws.Cells.Item(2,3).Value = "Ruhan, André, Pedro, Ivanilda, Silmara e Hélcio"
$ws.Cells.Item(2,4).Value = "Discussão acerca do problema principal referente a contração/aquisição e preenchimento do documento Business Case"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = "*Meta de redução ainda precisa de análise de dados históricos para coleta de médias. A partir destas será possível obter uma meta mais coerente."
$ws.Cells.Item(2,7).Value = "*Ruhan, Ivis , Andre - Criar planilha para organizar dados históricos das diversas categorias de contratações`n*Ivanilda , Silmara e Hélcio - Coleta de dados histórios desde 2012 das mais diversas categorias"
$ws.Cells.Item(2,8).Value = 42440
$ws.Cells.Item(2,8).NumberFormat = "mm/dd/yyyy"
$ws.Cells.Item(2,9).Value = "Em progresso"

$bodyRange = $ws.Range("A2:J2")
$bodyRange.Borders.LineStyle = 1
$bodyRange.Borders.Weight = 2
$bodyRange.WrapText = $true
$bodyRange.VerticalAlignment = -4160
$ws.Rows.Item(2).RowHeight = 105

$ws.Columns.Item(1).ColumnWidth = 10.29
$ws.Columns.Item(2).ColumnWidth = 17.86
$ws.Columns.Item(3).ColumnWidth = 17.86
$ws.Columns.Item(4).ColumnWidth = 25.86
$ws.Columns.Item(5).ColumnWidth = 10.29
$ws.Columns.Item(6).ColumnWidth = 29.43
$ws.Columns.Item(7).ColumnWidth = 31.71
$ws.Columns.Item(8).ColumnWidth = 9.71
$ws.Columns.Item(9).ColumnWidth = 12
$ws.Columns.Item(10).ColumnWidth = 27.43

$ws.Range("F16").Select()
